$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "CreatedAt: 2026-02-08T17:10:05"

$ws.Range("U4").Value = 368.62
$ws.Range("V4").Value = 358.26
$ws.Range("W4").Value = 206.45
$ws.Range("X4").Value = 212.02
$ws.Range("Y4").Value = 211.78
$ws.Range("Z4").Value = 298.85

$ws.Range("U6").Value = -26.91
$ws.Range("V6").Value = -25.08
$ws.Range("W6").Value = -12.59
$ws.Range("X6").Value = -12.3
$ws.Range("Y6").Value = -10.59
$ws.Range("Z6").Value = -13.15

$ws.Range("U9").Value = 343.34
$ws.Range("V9").Value = 337.45
$ws.Range("W9").Value = 195.4
$ws.Range("X9").Value = 208.08
$ws.Range("Y9").Value = 207.24
$ws.Range("Z9").Value = 292.69

$ws.Range("U11").Value = -52.19
$ws.Range("V11").Value = -45.89
$ws.Range("W11").Value = -23.64
$ws.Range("X11").Value = -16.23
$ws.Range("Y11").Value = -15.13
$ws.Range("Z11").Value = -19.32

$ws.Range("U14").Value = 68.86
$ws.Range("V14").Value = 62.84
$ws.Range("W14").Value = 195.4
$ws.Range("X14").Value = 208.08
$ws.Range("Y14").Value = 207.24
$ws.Range("Z14").Value = 292.69

$ws.Range("U15").Value = -274.48
$ws.Range("V15").Value = -274.61

$ws.Range("U16").Value = -52.19
$ws.Range("V16").Value = -45.89
$ws.Range("W16").Value = -23.64
$ws.Range("X16").Value = -16.23
$ws.Range("Y16").Value = -15.13
$ws.Range("Z16").Value = -19.32

$ws.Range("U19").Value = 182.9
$ws.Range("V19").Value = 104.38
$ws.Range("W19").Value = 142.04
$ws.Range("X19").Value = 114.3
$ws.Range("Y19").Value = 114.69
$ws.Range("Z19").Value = 77.90000000000001

$ws.Range("U20").Value = -179.97
$ws.Range("V20").Value = -248.6
$ws.Range("W20").Value = -61.34
$ws.Range("X20").Value = -94.95
$ws.Range("Y20").Value = -94.11
$ws.Range("Z20").Value = -216.73

$ws.Range("U21").Value = -32.66
$ws.Range("V21").Value = -30.36
$ws.Range("W21").Value = -15.66
$ws.Range("X21").Value = -15.07
$ws.Range("Y21").Value = -13.57
$ws.Range("Z21").Value = -17.38

$ws.Range("U24").Value = 362.87
$ws.Range("V24").Value = 352.98
$ws.Range("W24").Value = 203.38
$ws.Range("X24").Value = 209.25
$ws.Range("Y24").Value = 208.8
$ws.Range("Z24").Value = 294.62

$ws.Range("U26").Value = -32.66
$ws.Range("V26").Value = -30.36
$ws.Range("W26").Value = -15.66
$ws.Range("X26").Value = -15.07
$ws.Range("Y26").Value = -13.57
$ws.Range("Z26").Value = -17.38

$ws.Range("U29").Value = 177.33
$ws.Range("V29").Value = 98.94
$ws.Range("W29").Value = 138.88
$ws.Range("Y29").Value = 111.41
$ws.Range("Z29").Value = 72.97

$ws.Range("U30").Value = -179.97
$ws.Range("V30").Value = -248.6
$ws.Range("W30").Value = -61.34
$ws.Range("X30").Value = -94.95
$ws.Range("Y30").Value = -94.11
$ws.Range("Z30").Value = -216.73

$ws.Range("U31").Value = -38.23
$ws.Range("V31").Value = -35.8
$ws.Range("W31").Value = -18.82
$ws.Range("X31").Value = -18.14
$ws.Range("Y31").Value = -16.85
$ws.Range("Z31").Value = -22.31

$ws.Range("U34").Value = 61
$ws.Range("V34").Value = 57.86

$ws.Range("U35").Value = -274.48
$ws.Range("V35").Value = -274.61

$ws.Range("U36").Value = -60.05
$ws.Range("V36").Value = -50.87
$ws.Range("W36").Value = -26.73
$ws.Range("X36").Value = -16.42
$ws.Range("Y36").Value = -15.51
$ws.Range("Z36").Value = -19.32

$ws.Range("W37").Value = -191.31
$ws.Range("X37").Value = -206.89
$ws.Range("Y37").Value = -205.86
$ws.Range("Z37").Value = -291.69

$ws.Range("U39").Value = 368.62
$ws.Range("V39").Value = 358.26
$ws.Range("W39").Value = 206.45
$ws.Range("X39").Value = 212.02
$ws.Range("Y39").Value = 211.78
$ws.Range("Z39").Value = 298.85

$ws.Range("U41").Value = -26.91
$ws.Range("V41").Value = -25.08
$ws.Range("W41").Value = -12.59
$ws.Range("X41").Value = -12.3
$ws.Range("Y41").Value = -10.59
$ws.Range("Z41").Value = -13.15

$ws.Range("U44").Value = 391.23
$ws.Range("V44").Value = 379.17
$ws.Range("W44").Value = 216.02
$ws.Range("X44").Value = 220.78
$ws.Range("Y44").Value = 221.05
$ws.Range("Z44").Value = 311.38

$ws.Range("U46").Value = -4.3
$ws.Range("V46").Value = -4.17
$ws.Range("W46").Value = -3.02
$ws.Range("X46").Value = -3.53
$ws.Range("Y46").Value = -1.33
$ws.Range("Z46").Value = -0.62

$ws.Range("U49").Value = 372.09
$ws.Range("V49").Value = 358.59
$ws.Range("W49").Value = 214.33
$ws.Range("X49").Value = 232.93
$ws.Range("Y49").Value = 230.92
$ws.Range("Z49").Value = 325.68

$ws.Range("U51").Value = -23.44
$ws.Range("V51").Value = -24.74
$ws.Range("W51").Value = -4.72
$ws.Range("X51").Value = 8.619999999999999
$ws.Range("Y51").Value = 8.539999999999999
$ws.Range("Z51").Value = 13.68

$ws.Range("U54").Value = 381.05
$ws.Range("V54").Value = 371.81
$ws.Range("W54").Value = 212.87
$ws.Range("X54").Value = 220.57
$ws.Range("Y54").Value = 222.59
$ws.Range("Z54").Value = 309.84

$ws.Range("U56").Value = -14.48
$ws.Range("V56").Value = -11.53
$ws.Range("W56").Value = -6.17
$ws.Range("X56").Value = -3.75
$ws.Range("Y56").Value = 0.22
$ws.Range("Z56").Value = -2.17

$ws.Range("U59").Value = 404.01
$ws.Range("V59").Value = 391.16
$ws.Range("W59").Value = 221.93
$ws.Range("X59").Value = 227.04
$ws.Range("Y59").Value = 226.68
$ws.Range("Z59").Value = 319.02

$ws.Range("U61").Value = 8.48
$ws.Range("V61").Value = 7.82
$ws.Range("W61").Value = 2.89
$ws.Range("Y61").Value = 4.31
$ws.Range("Z61").Value = 7.02

$ws.Range("U64").Value = 412.87
$ws.Range("V64").Value = 400.56
$ws.Range("W64").Value = 226.52
$ws.Range("X64").Value = 231.49
$ws.Range("Y64").Value = 230.92
$ws.Range("Z64").Value = 324.33

$ws.Range("U66").Value = 17.34
$ws.Range("V66").Value = 17.22
$ws.Range("W66").Value = 7.48
$ws.Range("X66").Value = 7.18
$ws.Range("Y66").Value = 8.539999999999999
$ws.Range("Z66").Value = 12.32

$ws.Range("U69").Value = 410.3
$ws.Range("V69").Value = 398.07
$ws.Range("W69").Value = 224.43
$ws.Range("X69").Value = 229.36
$ws.Range("Y69").Value = 229.49
$ws.Range("Z69").Value = 323.32

$ws.Range("U71").Value = 14.77
$ws.Range("V71").Value = 14.73
$ws.Range("W71").Value = 5.39
$ws.Range("X71").Value = 5.05
$ws.Range("Y71").Value = 7.11
$ws.Range("Z71").Value = 11.32

$ws.Range("U74").Value = 405.26
$ws.Range("V74").Value = 392.36
$ws.Range("W74").Value = 223.28
$ws.Range("X74").Value = 228.19
$ws.Range("Y74").Value = 228.31
$ws.Range("Z74").Value = 320.99

$ws.Range("U76").Value = 9.73
$ws.Range("V76").Value = 9.02
$ws.Range("W76").Value = 4.24
$ws.Range("Y76").Value = 5.94
$ws.Range("Z76").Value = 8.99

$ws.Range("U79").Value = 395.53
$ws.Range("V79").Value = 383.34
$ws.Range("W79").Value = 219.04
$ws.Range("X79").Value = 224.32
$ws.Range("Y79").Value = 222.37
$ws.Range("Z79").Value = 312

$ws.Range("U84").Value = 373.49
$ws.Range("V84").Value = 366.48
$ws.Range("W84").Value = 210.01
$ws.Range("X84").Value = 218.63
$ws.Range("Y84").Value = 225.76
$ws.Range("Z84").Value = 304.1

$ws.Range("U86").Value = -22.04
$ws.Range("V86").Value = -16.86
$ws.Range("W86").Value = -9.029999999999999
$ws.Range("X86").Value = -5.68
$ws.Range("Y86").Value = 3.39
$ws.Range("Z86").Value = -7.91

$ws.Range("U89").Value = 357.3
$ws.Range("V89").Value = 347.54
$ws.Range("W89").Value = 200.22
$ws.Range("X89").Value = 206.17
$ws.Range("Y89").Value = 205.52
$ws.Range("Z89").Value = 289.7

$ws.Range("U91").Value = -38.23
$ws.Range("V91").Value = -35.8
$ws.Range("W91").Value = -18.82
$ws.Range("X91").Value = -18.14
$ws.Range("Y91").Value = -16.85
$ws.Range("Z91").Value = -22.31
